$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "62.809.21"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.96%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.974.81"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.13%  "
$r.Style = "Normal"

$r = $ws.Range("D4")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.07%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "596.55"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.21%  "
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "144.86"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.00%  "
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.973.37"
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.15%  "
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -0.16%  "
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "7.30"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +6.03%  "
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.80%  "
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.450"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.56%  "
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.0000237"
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +5.70%  "
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "33.61"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.59%  "
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.17%  "
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "3.471.34"
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.24%  "
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "62.686.58"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.83%  "
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.00%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.971.63"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.06%  "
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "443.18"
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.57%  "
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.10%  "
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.678"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.31%  "
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.34%  "
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "82.26"
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.07%  "
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "10.87"
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.83%  "
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "12.01"
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.15%  "
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.15"
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -1.82%  "
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.04%  "
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.11%  "
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "7.03"
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.88%  "
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -5.97%  "
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "26.55"
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.21%  "
$r.Style = "Normal"

$r = $ws.Range("D33")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.107"
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -0.98%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.16%  "
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.0₃0882"
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.92%  "
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -1.09%  "
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.60%  "
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "49.91"
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.54%  "
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -0.79%  "
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.02"
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +2.33%  "
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +1.10%  "
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -2.02%  "
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.281"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -0.70%  "
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "39.14"
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -4.19%  "
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "371.26"
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -1.43%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "2.702.69"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.21%  "
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "0.0341"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -1.24%  "
$r.Style = "Normal"

$r = $ws.Range("D48")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "134.08"
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.40%  "
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  +0.05%  "
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "23.22"
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -2.28%  "
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.Style = "Normal"
$r.NumberFormat = "@"
$r.Value = "  -0.51%  "
$r.Style = "Normal"

